# Updating test files to match the current format in beta
#
# This reproduces, via Excel COM interop, the semantic changes described by
# the target diff for the "optimization_parameters" sheet (and the resulting
# active-sheet/selection bookkeeping), namely:
#   - the stray duplicate "value" header cells in C1:F1 are cleared
#   - the "Model" label (A8) is renamed "production_function"
#   - a new "L_curve" / 0 row is inserted right after the production_function row
#   - the old "Deletion" / 0 / 3 row is removed entirely
#   - "optimization_parameters" becomes the active sheet/tab with C1:F1 selected
#
# (shared-string table renumbering/compaction and the dependent t="s" index
# shifts on every other sheet happen automatically as a consequence of these
# edits - the engine recomputes/saves the shared string table itself.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Row 1: drop the duplicated "value" header cells in C1:F1 ------------
$ws.Range("C1:F1").ClearContents()

# --- Row 8: "Model" -> "production_function" -----------------------------
$ws.Cells.Item(8, 1).Value = "production_function"

# --- Insert the new "L_curve" row directly below row 8 --------------------
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = "L_curve"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 2).NumberFormat = $ws.Cells.Item(2, 2).NumberFormat

# --- Remove the old "Deletion" row (now pushed down to row 17) ------------
$ws.Rows.Item(17).Delete()

# --- Make this sheet the active one, with C1:F1 selected ------------------
[void]$ws.Activate()
[void]$ws.Range("C1:F1").Select()
